# Updated legacy GSC export data
# The export dropped the oldest day (2025-10-11): its row is removed and
# every subsequent row shifts up by one. Because the refreshed export no
# longer has "Not indexed" / "Indexed" figures for the two oldest
# remaining days (2025-10-12 and 2025-10-13), those two cells are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the obsolete 2025-10-11 row; everything below shifts up one row.
$ws.Rows.Item(2).Delete()

# The two oldest remaining rows (now rows 2 and 3, for 2025-10-12 and
# 2025-10-13) no longer carry "Not indexed"/"Indexed" counts in the new
# export.
$ws.Range("B2:C3").Value = ""
